$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("蓝色光标", "白银有色", "利欧股份"),
    @("利欧股份", "利欧股份", "白银有色"),
    @("白银有色", "湖南白银", "平潭发展"),
    @("紫金矿业", "紫金矿业", "蓝色光标"),
    @("天奇股份", "蓝色光标", "华天科技"),
    @("湖南白银", "华天科技", "湖南白银"),
    @("华天科技", "湖南黄金", "协鑫集成"),
    @("航天电子", "铜陵有色", "紫金矿业"),
    @("湖南黄金", "航天电子", "泓博医药"),
    @("平潭发展", "浙文互联", "浙文互联"),
    @("达安基因", "平潭发展", "天奇股份"),
    @("浙文互联", "达安基因", "拓日新能"),
    @("协鑫集成", "协鑫集成", "航天发展"),
    @("三维通信", "三维通信", "航天电子"),
    @("铜陵有色", "拓日新能", "华胜天成"),
    @("拓日新能", "中国黄金", "四川黄金"),
    @("航天发展", "航天发展", "三维通信"),
    @("金风科技", "康强电子", "国晟科技"),
    @("优刻得-W", "天奇股份", "中国黄金"),
    @("中国黄金", "金风科技", "铜陵有色")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
